# Generate Report for Handoff
#
# A new handoff report run updates the "Latest Handoff Date/Datetime"
# values for the 0a38508a-bacc-493b-9519-d943db638858 file (row 5 on each
# sheet) on the Overview sheet and on each per-locale (zh-cn / de-de) sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-31-17 14:31:23"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-17 14:31:20"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-17 14:31:23"
